$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.348.96'
$ws.Range("E2").Value = '  -0.73%  '

$ws.Range("D3").Value = '1.713.80'
$ws.Range("E3").Value = '  -1.36%  '

$ws.Range("E4").Value = '  -0.36%  '

$ws.Range("D5").Value = '240.23'
$ws.Range("E5").Value = '  -2.70%  '

$ws.Range("D6").Value = '0.9967'
$ws.Range("E6").Value = '  -0.35%  '

$ws.Range("D7").Value = '0.4866'
$ws.Range("E7").Value = '  -1.52%  '

$ws.Range("D8").Value = '0.2584'
$ws.Range("E8").Value = '  -3.14%  '

$ws.Range("D9").Value = '0.06175'
$ws.Range("E9").Value = '  -1.90%  '

$ws.Range("D10").Value = '1.712.38'
$ws.Range("E10").Value = '  -1.11%  '

$ws.Range("D11").Value = '0.06946'
$ws.Range("E11").Value = '  -1.41%  '

$ws.Range("D12").Value = '15.49'
$ws.Range("E12").Value = '  -1.45%  '

$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = '0.5972'
$ws.Range("E13").Value = '  -2.50%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '4.468'
$ws.Range("E14").Value = '  -2.78%  '

$ws.Range("D15").Value = '76.42'
$ws.Range("E15").Value = '  -1.34%  '

$ws.Range("D16").Value = '0.9965'
$ws.Range("E16").Value = '  -0.37%  '

$ws.Range("D17").Value = '26.246.21'
$ws.Range("E17").Value = '  -1.10%  '

$ws.Range("D18").Value = '0.9961'
$ws.Range("E18").Value = '  -0.37%  '

$ws.Range("D19").Value = '0.000007087'
$ws.Range("E19").Value = '  -3.48%  '

$ws.Range("D20").Value = '11.25'
$ws.Range("E20").Value = '  -2.56%  '

$ws.Range("D21").Value = '1.933.20'
$ws.Range("E21").Value = '  -1.19%  '

$ws.Range("D22").Value = '4.409'
$ws.Range("E22").Value = '  -3.98%  '

$ws.Range("D23").Value = '8.432'
$ws.Range("E23").Value = '  -3.11%  '

$ws.Range("D24").Value = '5.047'
$ws.Range("E24").Value = '  -3.90%  '

$ws.Range("D25").Value = '136.35'
$ws.Range("E25").Value = '  -2.67%  '

$ws.Range("E26").Value = '  -2.08%  '

$ws.Range("D27").Value = '1.394'
$ws.Range("E27").Value = '  -1.85%  '

$ws.Range("E28").Value = '  -1.88%  '

$ws.Range("D29").Value = '105.37'
$ws.Range("E29").Value = '  -2.51%  '

$ws.Range("D30").Value = '3.876'
$ws.Range("E30").Value = '  -4.20%  '

$ws.Range("D31").Value = '0.07931'
$ws.Range("E31").Value = '  -1.83%  '

$ws.Range("D32").Value = '3.607'
$ws.Range("E32").Value = '  -2.95%  '

$ws.Range("D33").Value = '0.04427'
$ws.Range("E33").Value = '  -3.64%  '

$ws.Range("B34").Value = 'Frax'
$ws.Range("C34").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D34").Value = '0.9957'
$ws.Range("E34").Value = '  -0.42%  '

$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '2.596'
$ws.Range("E35").Value = '  -0.58%  '

$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '0.9900'
$ws.Range("E36").Value = '  -1.85%  '

$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '0.6164'
$ws.Range("E37").Value = '  -3.10%  '

$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '0.9342'
$ws.Range("E38").Value = '  +4.01%  '

$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '1.977'
$ws.Range("E39").Value = '  -1.99%  '

$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '2.372'
$ws.Range("E40").Value = '  -1.18%  '

$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").Value = '0.9955'
$ws.Range("E41").Value = '  -0.91%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '0.01471'

$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '99.69'
$ws.Range("E43").Value = '  -2.22%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '5.373'
$ws.Range("E44").Value = '  -0.55%  '

$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '0.3800'
$ws.Range("E45").Value = '  -2.76%  '

$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '6.813'
$ws.Range("E46").Value = '  -0.85%  '

$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '0.1149'
$ws.Range("E47").Value = '  -3.30%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.05343'
$ws.Range("E48").Value = '  -1.01%  '

$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '30.55'
$ws.Range("E49").Value = '  -0.05%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '7.677'
$ws.Range("E50").Value = '  -1.82%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '51.06'
$ws.Range("E51").Value = '  -1.46%  '
